$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# capture existing comment texts before mutating
$texts = @{}
foreach ($addr in @("F11","G11","H11","I11","J11","K11","L11","M11","N11","O11")) {
  $texts[$addr] = $ws.Range($addr).Comment.Text()
}

# delete E11 comment (GENOMIC SOURCE one)
$ws.Range("E11").Comment.Delete()

# delete F11:O11 comments
foreach ($addr in @("F11","G11","H11","I11","J11","K11","L11","M11","N11","O11")) {
  $ws.Range($addr).Comment.Delete()
}

# add comments shifted left by one: F11 -> E11, etc.
$mapping = @{
  "E11" = $texts["F11"];
  "F11" = $texts["G11"];
  "G11" = $texts["H11"];
  "H11" = $texts["I11"];
  "I11" = $texts["J11"];
  "J11" = $texts["K11"];
  "K11" = $texts["L11"];
  "L11" = "Path to a file on a users computer";
  "M11" = $texts["M11"];
  "N11" = $texts["O11"];
}
foreach ($addr in $mapping.Keys) {
  $ws.Range($addr).AddComment($mapping[$addr])
}

Write-Host $ws.Comments.Count
foreach ($cm in $ws.Comments) {
  Write-Host ($cm.Parent.Address() + " => " + $cm.Text())
}
